$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '50.916.75'
$ws.Range('E2').Value = '  -1.69%  '
$ws.Range('D3').Value = '2.906.85'
$ws.Range('E3').Value = '  -1.85%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '369.18'
$ws.Range('E5').Value = '  +4.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.63'
$ws.Range('E6').Value = '  -5.50%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.540'
$ws.Range('E7').Value = '  -2.82%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.581'
$ws.Range('E9').Value = '  -4.39%  '
$ws.Range('E10').Value = '  -3.40%  '
$ws.Range('E11').Value = '  +0.49%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0832'
$ws.Range('E12').Value = '  -2.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.22'
$ws.Range('E13').Value = '  -4.22%  '
$ws.Range('D14').Value = '3.357.07'
$ws.Range('E14').Value = '  -1.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.34'
$ws.Range('E15').Value = '  -3.68%  '
$ws.Range('D16').Value = '2.900.03'
$ws.Range('E16').Value = '  -2.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.922'
$ws.Range('E17').Value = '  -6.19%  '
$ws.Range('D18').Value = '50.842.51'
$ws.Range('E18').Value = '  -1.77%  '
$ws.Range('E19').Value = '  -5.87%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.16'
$ws.Range('E20').Value = '  -3.45%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.85'
$ws.Range('E21').Value = '  -4.49%  '
$ws.Range('D22').Value = '0.0₃0939'
$ws.Range('E22').Value = '  -2.99%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.87'
$ws.Range('E23').Value = '  -1.83%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '257.72'
$ws.Range('E24').Value = '  -1.57%  '
$ws.Range('E25').Value = '  -1.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.17'
$ws.Range('E26').Value = '  -2.87%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.165'
$ws.Range('E28').Value = '  -4.90%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '25.52'
$ws.Range('E29').Value = '  -4.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.08'
$ws.Range('E30').Value = '  -3.26%  '
$ws.Range('E31').Value = '  -4.51%  '
$ws.Range('E32').Value = '  +2.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '9.85'
$ws.Range('E33').Value = '  -4.01%  '
$ws.Range('E34').Value = '  -2.28%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '51.31'
$ws.Range('E35').Value = '  +1.49%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '34.05'
$ws.Range('E36').Value = '  -4.83%  '
$ws.Range('E37').Value = '  +0.65%  '
$ws.Range('E38').Value = '  -3.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.97'
$ws.Range('E39').Value = '  -6.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '16.99'
$ws.Range('E40').Value = '  -4.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.58'
$ws.Range('E41').Value = '  -3.24%  '
$ws.Range('E42').Value = '  -6.11%  '
$ws.Range('E43').Value = '  -3.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '119.06'
$ws.Range('E44').Value = '  -3.49%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.88'
$ws.Range('E45').Value = '  -1.72%  '
$ws.Range('E46').Value = '  -1.41%  '
$ws.Range('D47').Value = '2.016.42'
$ws.Range('E47').Value = '  -4.32%  '
$ws.Range('E48').Value = '  +0.14%  '
$ws.Range('E49').Value = '  -6.07%  '
$ws.Range('D50').Value = '3.191.35'
$ws.Range('E50').Value = '  -1.13%  '
$ws.Range('E51').Value = '  -0.69%  '
